# feat: add 2022-Q1 data
#
# Before: 3 sheets -> "2021-Q1", "2021-Q4", "总计" (summary table).
# After:  4 sheets -> "2021-Q1", "2021-Q4", "2022-Q1" (new fund-holdings
#         detail sheet, reusing the old "总计" sheet's slot/id), "总计"
#         (brand new summary sheet, appended right after "2022-Q1", with
#         an extra row for the new quarter).

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet into the new "2022-Q1"
# detail sheet (keeps the same sheetId/slot - just renamed + refilled).
# -----------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# -----------------------------------------------------------------------
# Step 2: while $q1 still holds the old "总计" table (and its formatting
# - identical sheetPr/pageMargins to the other sheets already), copy it
# to create the brand-new "总计" summary sheet right after it. Grab the
# new sheet by position (Index+1) rather than by its auto-generated
# "... (2)" name, so this doesn't depend on naming/locale conventions.
# -----------------------------------------------------------------------
$q1Index = $q1.Index
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item($q1Index + 1)
$total.Name = "总计"

# --- Populate the new "总计" summary sheet ------------------------------
# Insert a fresh data row on top for the new quarter; existing rows
# (2021-Q4, 2021-Q1) shift down automatically.
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(3, 1).Copy($total.Cells.Item(2, 1))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.06

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2

# -----------------------------------------------------------------------
# Step 3: refill "2022-Q1" with the new fund-holdings detail table
# (it still contains the stale "总计" rows at this point).
# -----------------------------------------------------------------------
$q1.Rows.Item(3).Delete()
$q1.Rows.Item(2).Delete()

# Header row, styled like every other detail sheet (copy from "2021-Q4").
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("B1:H1").Copy($q1.Range("B1:H1"))

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Data row 2 (style copied from "2021-Q4" row 2: A column centered/bordered,
# the rest plain - matches the rest of the workbook's detail sheets).
$src.Range("A2:H2").Copy($q1.Range("A2:H2"))

$q1.Cells.Item(2, 1).Value = 0

$q1.Range("B2").NumberFormat = "@"
$q1.Cells.Item(2, 2).Value = "160620"

$q1.Cells.Item(2, 3).Value = "鹏华中证A股资源产业指数（LOF）"

$q1.Range("D2:G2").NumberFormat = "@"
$q1.Cells.Item(2, 4).Value = "2.77"
$q1.Cells.Item(2, 5).Value = "94.14"
$q1.Cells.Item(2, 6).Value = "2.33"
$q1.Cells.Item(2, 7).Value = "0.0645"

$q1.Cells.Item(2, 8).Value = 10
